# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates DAMSLTag (column I) and DialogAct (column J) values for a set of rows
# on Sheet1 to reflect the re-annotated dialog acts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 11;  DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 37;  DAMSLTag = "b";  DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 47;  DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 66;  DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 70;  DAMSLTag = "ba"; DialogAct = "Appreciation" },
    @{ Row = 72;  DAMSLTag = "qy"; DialogAct = "Yes-No-Question" },
    @{ Row = 114; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 120; DAMSLTag = "%";  DialogAct = "Uninterpretable" },
    @{ Row = 136; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 137; DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 138; DAMSLTag = "%";  DialogAct = "Uninterpretable" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

$wb.Save()
